$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data rows (A2:B23) descending by the value column (B), keeping
# the header row (row 1) untouched.
$dataRange = $ws.Range("A2:B23")
$sortKey = $ws.Range("B2:B23")
$dataRange.Sort($sortKey, 2)

# After sorting descending by value, the two lowest-value languages
# (Swedish and Uzbek) land in rows 22 and 23 - remove them entirely.
$ws.Rows("22:23").Delete()
